# Detailed Engineering.xlsx -- "tip vs. slip" analysis update
#
# 1. "joints and bearings" sheet: mark the two safety-factor helper cells
#    (C90 / C131) with an "Overdriving" annotation, and drop the feet-tube
#    outer diameter (B127) from 15 mm to 13 mm (cascades into B129/B130/
#    B131/B133/B134/B135/B109/B110 via the existing formulas).
# 2. Add a brand-new "tip vs. slip" worksheet at the end of the workbook
#    with the tip-vs-slip overdriving calculation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update "joints and bearings"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("joints and bearings")

$ws2.Range("C90").Value = "Overdriving"
$ws2.Range("C131").Value = "Overdriving"

$ws2.Range("B127").Value = 13

$ws2.Activate()
$ws2.Range("D105").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Add the "tip vs. slip" worksheet as the last tab
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "tip vs. slip"

$ws3.Range("A1").Value = "Coefficient of friction, wood on wood"
$ws3.Range("B1").Value = 0.5

$ws3.Range("A2").Value = "CG Height (mm)"
$ws3.Range("B2").Value = 300
$ws3.Range("C2").Value = "approximating as midway up to first order"

$ws3.Range("A3").Value = "Mass of entire structure (kg)"
$ws3.Range("B3").Value = 8.78

$ws3.Range("A4").Value = "Dead weight of structure (N)"
$ws3.Range("B4").Formula = "=B3*9.81"

$ws3.Range("A5").Value = "Desktop height at top of travel (mm)"

$ws3.Columns.Item(1).ColumnWidth = 34.09

$ws3.Range("A6").Select() | Out-Null
$ws3.Activate()
